$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy style from an existing header cell (e.g. E1) to the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean FALSE values for rows 2-7, columns F:H
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}
